# Hortaliza, Macroferia Regional de Talca - Poroto verde
# Weekly update: a new price record is inserted at row 129 (pushing the
# existing rows 129-155 down to 130-156).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 129; this shifts rows 129:155
# down to 130:156 and extends the used range to A1:R156.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with this week's record.
$ws.Cells.Item(129, 1).Value = 5
$ws.Cells.Item(129, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(129, 3).Value = "Maule"
$ws.Cells.Item(129, 4).Value = 44644
$ws.Cells.Item(129, 5).Value = 7
$ws.Cells.Item(129, 6).Value = 100112031
$ws.Cells.Item(129, 7).Value = "Poroto verde"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 150
$ws.Cells.Item(129, 11).Value = 25000
$ws.Cells.Item(129, 12).Value = 25000
$ws.Cells.Item(129, 13).Value = 25000
$ws.Cells.Item(129, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(129, 15).Value = "Región del Maule"
$ws.Cells.Item(129, 16).Value = 1000
$ws.Cells.Item(129, 17).Value = 25
$ws.Cells.Item(129, 18).Value = "Hortaliza"
